# Update the "Förändrad" (Changed) date column (C) for every data row
# from 2023-09-01 (serial 45170) to 2023-09-05 (serial 45174).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1
if ($lastRow -lt 2) { $lastRow = 108 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45170) {
        $cell.Value = 45174
    }
}
